$wb = $excel.ActiveWorkbook

# Clone the most recent weekly sheet (inherits header styles/col widths/page
# setup) and place the copy at the end of the tab strip
$srcSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcSheet.Copy($null, $srcSheet)

# The copy is now the last sheet; grab + rename it
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "magapoke_2025-11-12"

# Header row (kept identical to the source sheet, rewritten explicitly for clarity)
$ws.Cells.Item(1,1).Value = "rank"
$ws.Cells.Item(1,2).Value = "title"

# Ranked rows 2..101 -> rank 1..100 with this week's titles
$titles = @(
    'WIND BREAKER',
    '信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐＆『ざまぁ！』します！',
    '東京卍リベンジャーズ',
    'ギルティサークル',
    '宇宙兄弟',
    'みいちゃんと山田さん',
    '島耕作',
    '薫る花は凛と咲く',
    'FAIRY TAIL 100 YEARS QUEST',
    'イレギュラーズ',
    '転生したら第七王子だったので、気ままに魔術を極めます',
    '愛妻の裏アカ',
    '君が僕らを悪魔と呼んだ頃',
    '十字架のろくにん',
    '黄昏町プリズナーズ',
    '南海トラフ巨大地震',
    'ハードワーカー中田',
    '魔術ギルド総帥～生まれ変わって今更やり直す2度目の学院生活～',
    'ドラハチ',
    '味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す',
    'アルキメデスの大戦',
    'K-9~警視庁公安部公安第9課異能対策係~',
    'ひゃくえむ。',
    '転生貴族、鑑定スキルで成り上がる～弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた～',
    '蒼く染めろ',
    '魔女と傭兵',
    '辺境の薬師、都でSランク冒険者となる～英雄村の少年がチート薬で無自覚無双〜',
    '触手魔術師の成り上がり',
    'さわらないで小手指くん',
    'ジュミドロ',
    'グラぱらっ！',
    '食糧人類-Starving Anonymous-',
    'おやすみ ふみさん',
    'ともだちづくり',
    '追放された転生王子、『自動製作《オートクラフト》』スキルで領地を爆速で開拓し最強の村を作ってしまう〜最強クラフトスキルで始める、楽々領地開拓スローライフ〜',
    'せいぶつ部の田辺くん',
    '阿武ノーマル',
    '不遇職【鑑定士】が実は最強だった～奈落で鍛えた最強の【神眼】で無双する～',
    'ナキナギ',
    '追放されなかった男　～二度目の人生は土下座から始まりました～',
    'ハナバス　苔石花江のバスケ論',
    'いじめるヤバイ奴',
    'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！',
    '異世界ウォーキング',
    '幼馴染とはラブコメにならない',
    '時々ボソッとロシア語でデレる隣のアーリャさん',
    'Aランクパーティを離脱した俺は、元教え子たちと迷宮深部を目指す。',
    '不遇職『鍛冶師』だけど最強です ～気づけば何でも作れるようになっていた男ののんびりスローライフ～',
    'Destiny Unchain Online 〜吸血鬼少女となって、やがて『赤の魔王』と呼ばれるようになりました〜',
    '東京卍リベンジャーズ～場地圭介からの手紙～',
    'この世界がいずれ滅ぶことを、俺だけが知っている～モンスターが現れた世界で、死に戻りレベルアップ～',
    '五輪の女神さま 〜なでしこ寮のメダルごはん〜',
    '鉱石令嬢〜没落した悪役令嬢が炭鉱で一山当てるまでのお話〜',
    'デッドアカウント',
    '念願の悪役令嬢（ラスボス）の身体を手に入れたぞ！',
    'ストーカー行為がバレて人生終了男',
    'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜',
    '春くらり',
    '恋ニ非ズ',
    '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～',
    'お嬢様の僕',
    '可愛いだけじゃない式守さん',
    'アオバノバスケ',
    'となりの黒川さん',
    '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～',
    '四刀流の最強配信者～やり込んだVRゲームの設定が現実世界に反映されたので、廃止予定だった戦闘職で無双します～',
    '降り積もれ孤独な死よ',
    'ブルーロック',
    'ダメスキル【自動機能】が覚醒しました～あれ、ギルドのスカウトの皆さん、俺を「いらない」って言ってませんでした？～',
    '田んぼで拾った女騎士、田舎で俺の嫁だと思われている',
    '屋根の下のアルテミス',
    '東京ネオンスキャンダル',
    '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～',
    'デスティニーラバーズ',
    'リスナーに騙されてダンジョンの最下層から脱出RTAすることになった',
    'なれの果ての僕ら',
    'シャングリラ・フロンティア～クソゲーハンター、神ゲーに挑まんとす～',
    '黒猫と魔女の教室',
    '最弱な僕は＜壁抜けバグ＞で成り上がる～壁をすり抜けたら、初回クリア報酬を無限回収できました！～',
    '勇者と呼ばれた後に　―そして無双男は家族を創る―',
    'それがメイドのカンナです',
    '復讐の教科書',
    'GALAXIAS',
    'ヒロインは絶望しました。',
    '冰剣の魔術師が世界を統べる〜世界最強の魔術師である少年は、魔術学院に入学する〜',
    '我間乱 ―修羅―',
    'イジらないで、長瀞さん',
    'DAYS外伝',
    'MYS',
    'はっちぽっちぱんち',
    'インフェクション',
    '魁の花巫女',
    '邪魔な初級職を追放したら、大変なことになっちゃったんですけど！？～追放された初級職【アイテム師】が自分の居場所を見つけるまで外伝～',
    'はじめの一歩',
    '剣帝学院の魔眼賢者',
    '母という呪縛 娘という牢獄',
    '彼女、お借りします',
    '〈小市民〉 春期限定いちごタルト事件',
    '人間消失',
    'ハプスブルク家の華麗なる受難'
)

for ($i = 0; $i -lt $titles.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $titles[$i]
}

# Copying a sheet makes the new copy active; restore the original active tab
# (the first sheet) so the workbook-level selection state is unchanged
$wb.Worksheets.Item(1).Select()
